$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: opening line of the notice.
#   "In the " -> "In the" / " county court at" / " "  (three runs in the
#   canonical XML, all sharing the same bold Arial 12pt formatting as the
#   original run, immediately followed by the existing "<<hearingSiteName>>"
#   merge field run).
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("In the ", $false, $false, $false, $false, $false, $true, 1, $false, `
    "In the county court at ", 2)
Write-Host "Change1 (opening 'In the' line): $found1"

# ---------------------------------------------------------------------------
# Change 2: hearing-location sentence.
#   " of the claimant's claim will take place " + "at"
#     -> " of the claimant's claim will take place" + " in the county court sitting at"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("of the claimant" + [char]8217 + "s claim will take place at", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "of the claimant" + [char]8217 + "s claim will take place in the county court sitting at", 2)
Write-Host "Change2 (hearing location sentence): $found2"

# ---------------------------------------------------------------------------
# Change 3: drop the stale <w:lastRenderedPageBreak/> marker in front of the
# video/telephone hearing guidance paragraph. Re-writing the run's text
# (identity replace) forces the host to rebuild the run, which drops the
# obsolete lastRenderedPageBreak marker while leaving the text untouched.
# ---------------------------------------------------------------------------
$guidanceText = "Please use the following URL link to review the video (Cloud Video Platform) and telephone (BT meet me) hearing guidance if required."
$rng3 = $d.Content
$found3 = $rng3.Find.Execute($guidanceText, $false, $false, $false, $false, $false, $true, 1, $false, `
    $guidanceText, 2)
Write-Host "Change3 (drop stale lastRenderedPageBreak): $found3"
